$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il16"
$ws.Cells.Item(2, 3).Value = "Grin2d"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 6.566789
$ws.Cells.Item(2, 8).Value = 19.700367
$ws.Cells.Item(2, 9).Value = 0.3634074580963036
$ws.Cells.Item(2, 10).Value = 0.3634074580963036
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06287633333333333
$ws.Cells.Item(2, 14).Value = 0.188629
$ws.Cells.Item(2, 15).Value = 0.01674921743188729
$ws.Cells.Item(2, 16).Value = 0.01674921743188729
$ws.Cells.Item(2, 17).Value = 0.4128956140936666
$ws.Cells.Item(2, 18).Value = 3.716060526843
$ws.Cells.Item(2, 19).Value = 0.006086790532024458
$ws.Cells.Item(2, 20).Value = 0.006086790532024459

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il16"
$ws.Cells.Item(3, 3).Value = "Grin2d"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 6.566789
$ws.Cells.Item(3, 8).Value = 19.700367
$ws.Cells.Item(3, 9).Value = 0.3634074580963036
$ws.Cells.Item(3, 10).Value = 0.3634074580963036
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 3.112268
$ws.Cells.Item(3, 14).Value = 9.336804
$ws.Cells.Item(3, 15).Value = 0.8290568275022134
$ws.Cells.Item(3, 16).Value = 0.8290568275022134
$ws.Cells.Item(3, 17).Value = 20.437607267452
$ws.Cells.Item(3, 18).Value = 183.938465407068
$ws.Cells.Item(3, 19).Value = 0.301285434299965
$ws.Cells.Item(3, 20).Value = 0.301285434299965

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il16"
$ws.Cells.Item(4, 3).Value = "Grin2d"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 6.566789
$ws.Cells.Item(4, 8).Value = 19.700367
$ws.Cells.Item(4, 9).Value = 0.3634074580963036
$ws.Cells.Item(4, 10).Value = 0.3634074580963036
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 0.578842
$ws.Cells.Item(4, 14).Value = 1.736526
$ws.Cells.Item(4, 15).Value = 0.1541939550658993
$ws.Cells.Item(4, 16).Value = 0.1541939550658993
$ws.Cells.Item(4, 17).Value = 3.801133278338
$ws.Cells.Item(4, 18).Value = 34.210199505042
$ws.Cells.Item(4, 19).Value = 0.0560352332643141
$ws.Cells.Item(4, 20).Value = 0.05603523326431411

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il16"
$ws.Cells.Item(5, 3).Value = "Grin2d"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 6.440526333333334
$ws.Cells.Item(5, 8).Value = 19.321579
$ws.Cells.Item(5, 9).Value = 0.3564200560729107
$ws.Cells.Item(5, 10).Value = 0.3564200560729107
$ws.Cells.Item(5, 11).Value = 1.0
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.06287633333333333
$ws.Cells.Item(5, 14).Value = 0.188629
$ws.Cells.Item(5, 15).Value = 0.01674921743188729
$ws.Cells.Item(5, 16).Value = 0.01674921743188729
$ws.Cells.Item(5, 17).Value = 0.4049566805767777
$ws.Cells.Item(5, 18).Value = 3.644610125191
$ws.Cells.Item(5, 19).Value = 0.005969757016250641
$ws.Cells.Item(5, 20).Value = 0.005969757016250641

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il16"
$ws.Cells.Item(6, 3).Value = "Grin2d"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 6.440526333333334
$ws.Cells.Item(6, 8).Value = 19.321579
$ws.Cells.Item(6, 9).Value = 0.3564200560729107
$ws.Cells.Item(6, 10).Value = 0.3564200560729107
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 3.112268
$ws.Cells.Item(6, 14).Value = 9.336804
$ws.Cells.Item(6, 15).Value = 0.8290568275022134
$ws.Cells.Item(6, 16).Value = 0.8290568275022134
$ws.Cells.Item(6, 17).Value = 20.04464401039067
$ws.Cells.Item(6, 18).Value = 180.401796093516
$ws.Cells.Item(6, 19).Value = 0.2954924809459684
$ws.Cells.Item(6, 20).Value = 0.2954924809459684

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il16"
$ws.Cells.Item(7, 3).Value = "Grin2d"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 6.440526333333334
$ws.Cells.Item(7, 8).Value = 19.321579
$ws.Cells.Item(7, 9).Value = 0.3564200560729107
$ws.Cells.Item(7, 10).Value = 0.3564200560729107
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 0.578842
$ws.Cells.Item(7, 14).Value = 1.736526
$ws.Cells.Item(7, 15).Value = 0.1541939550658993
$ws.Cells.Item(7, 16).Value = 0.1541939550658993
$ws.Cells.Item(7, 17).Value = 3.728047143839333
$ws.Cells.Item(7, 18).Value = 33.552424294554
$ws.Cells.Item(7, 19).Value = 0.05495781811069169
$ws.Cells.Item(7, 20).Value = 0.05495781811069169

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Il16"
$ws.Cells.Item(8, 3).Value = "Grin2d"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 4.846355333333333
$ws.Cells.Item(8, 8).Value = 14.539066
$ws.Cells.Item(8, 9).Value = 0.2681983040292799
$ws.Cells.Item(8, 10).Value = 0.2681983040292799
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06287633333333333
$ws.Cells.Item(8, 14).Value = 0.188629
$ws.Cells.Item(8, 15).Value = 0.01674921743188729
$ws.Cells.Item(8, 16).Value = 0.01674921743188729
$ws.Cells.Item(8, 17).Value = 0.3047210533904444
$ws.Cells.Item(8, 18).Value = 2.742489480514
$ws.Cells.Item(8, 19).Value = 0.004492111709049821
$ws.Cells.Item(8, 20).Value = 0.004492111709049822

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Il16"
$ws.Cells.Item(9, 3).Value = "Grin2d"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 4.846355333333333
$ws.Cells.Item(9, 8).Value = 14.539066
$ws.Cells.Item(9, 9).Value = 0.2681983040292799
$ws.Cells.Item(9, 10).Value = 0.2681983040292799
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 3.112268
$ws.Cells.Item(9, 14).Value = 9.336804
$ws.Cells.Item(9, 15).Value = 0.8290568275022134
$ws.Cells.Item(9, 16).Value = 0.8290568275022134
$ws.Cells.Item(9, 17).Value = 15.08315662056267
$ws.Cells.Item(9, 18).Value = 135.748409585064
$ws.Cells.Item(9, 19).Value = 0.2223516350799889
$ws.Cells.Item(9, 20).Value = 0.2223516350799889

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Il16"
$ws.Cells.Item(10, 3).Value = "Grin2d"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 4.846355333333333
$ws.Cells.Item(10, 8).Value = 14.539066
$ws.Cells.Item(10, 9).Value = 0.2681983040292799
$ws.Cells.Item(10, 10).Value = 0.2681983040292799
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.578842
$ws.Cells.Item(10, 14).Value = 1.736526
$ws.Cells.Item(10, 15).Value = 0.1541939550658993
$ws.Cells.Item(10, 16).Value = 0.1541939550658993
$ws.Cells.Item(10, 17).Value = 2.805274013857333
$ws.Cells.Item(10, 18).Value = 25.247466124716
$ws.Cells.Item(10, 19).Value = 0.04135455724024117
$ws.Cells.Item(10, 20).Value = 0.04135455724024117

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Il16"
$ws.Cells.Item(11, 3).Value = "Grin2d"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 0.216374
$ws.Cells.Item(11, 8).Value = 0.6491220000000001
$ws.Cells.Item(11, 9).Value = 0.01197418180150597
$ws.Cells.Item(11, 10).Value = 0.01197418180150597
$ws.Cells.Item(11, 11).Value = 1.0
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.06287633333333333
$ws.Cells.Item(11, 14).Value = 0.188629
$ws.Cells.Item(11, 15).Value = 0.01674921743188729
$ws.Cells.Item(11, 16).Value = 0.01674921743188729
$ws.Cells.Item(11, 17).Value = 0.01360480374866667
$ws.Cells.Item(11, 18).Value = 0.122443233738
$ws.Cells.Item(11, 19).Value = 0.0002005581745623714
$ws.Cells.Item(11, 20).Value = 0.0002005581745623714

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Il16"
$ws.Cells.Item(12, 3).Value = "Grin2d"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 0.216374
$ws.Cells.Item(12, 8).Value = 0.6491220000000001
$ws.Cells.Item(12, 9).Value = 0.01197418180150597
$ws.Cells.Item(12, 10).Value = 0.01197418180150597
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 3.112268
$ws.Cells.Item(12, 14).Value = 9.336804
$ws.Cells.Item(12, 15).Value = 0.8290568275022134
$ws.Cells.Item(12, 16).Value = 0.8290568275022134
$ws.Cells.Item(12, 17).Value = 0.6734138762320002
$ws.Cells.Item(12, 18).Value = 6.060724886088002
$ws.Cells.Item(12, 19).Value = 0.00992727717629128
$ws.Cells.Item(12, 20).Value = 0.00992727717629128

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Il16"
$ws.Cells.Item(13, 3).Value = "Grin2d"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 0.216374
$ws.Cells.Item(13, 8).Value = 0.6491220000000001
$ws.Cells.Item(13, 9).Value = 0.01197418180150597
$ws.Cells.Item(13, 10).Value = 0.01197418180150597
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 0.578842
$ws.Cells.Item(13, 14).Value = 1.736526
$ws.Cells.Item(13, 15).Value = 0.1541939550658993
$ws.Cells.Item(13, 16).Value = 0.1541939550658993
$ws.Cells.Item(13, 17).Value = 0.125246358908
$ws.Cells.Item(13, 18).Value = 1.127217230172
$ws.Cells.Item(13, 19).Value = 0.001846346450652321
$ws.Cells.Item(13, 20).Value = 0.001846346450652321

Write-Host "Edit complete"